$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "33.947.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.777.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.554"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.289"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0699"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.034.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.764.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.936.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0781"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  +3.24%  "
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.390.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.659"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.50%  "
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.92%  "
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.910"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "77.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₆0146"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +22.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +13.19%  "
$ws.Range("E46").Value = "  +4.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "108.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.51%  "
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.934.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("E51").Value = "  +0.54%  "
